$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0.2151403427124023
$ws.Range("C1").Value = 0.2862262725830078
$ws.Range("D1").Value = 0.3079617023468018
$ws.Range("E1").Value = 0.3019604682922363
$ws.Range("F1").Value = 0.1795210838317871
$ws.Range("G1").Value = 0.2303841114044189
$ws.Range("H1").Value = -0.1047201156616211
$ws.Range("I1").Value = -0.07293176651000977
$ws.Range("J1").Value = -0.2441329956054688
$ws.Range("K1").Value = -0.05684876441955566
$ws.Range("L1").Value = -0.1795210838317871
$ws.Range("M1").Value = 0.6895039081573486
$ws.Range("N1").Value = 0.7017560005187988
$ws.Range("O1").Value = 0.5393245220184326
$ws.Range("P1").Value = 0.8630392551422119
